$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list refresh: update Price (D) and Volume(1h) (E) columns,
# and fix two mis-ordered rows (EthereumClassic/ImmutableX, dogwifhat/InjectiveProtocol).

$ws.Range('D2').Value = '69.336.39'
$ws.Range('E2').Value = '  -0.04%  '
$ws.Range('D3').Value = '3.684.27'
$ws.Range('E3').Value = '  -0.01%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').Value = "'678.23"
$ws.Range('E5').Value = '  -1.25%  '
$ws.Range('D6').Value = "'159.05"
$ws.Range('E6').Value = '  -2.22%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('E8').Value = '  -1.05%  '
$ws.Range('E9').Value = '  -1.11%  '
$ws.Range('D10').Value = "'7.06"
$ws.Range('E10').Value = '  -4.32%  '
$ws.Range('E11').Value = '  -1.92%  '
$ws.Range('E12').Value = '  -3.29%  '
$ws.Range('D13').Value = '4.308.28'
$ws.Range('E13').Value = '  +0.07%  '
$ws.Range('D14').Value = "'32.39"
$ws.Range('E14').Value = '  -3.52%  '
$ws.Range('D15').Value = '3.675.64'
$ws.Range('E15').Value = '  -0.06%  '
$ws.Range('D16').Value = '69.294.24'
$ws.Range('E16').Value = '  -0.19%  '
$ws.Range('E17').Value = '  +1.88%  '
$ws.Range('D18').Value = "'16.04"
$ws.Range('E18').Value = '  -1.88%  '
$ws.Range('D19').Value = "'6.42"
$ws.Range('E19').Value = '  -3.17%  '
$ws.Range('D20').Value = "'468.50"
$ws.Range('E20').Value = '  -3.26%  '
$ws.Range('D21').Value = "'9.97"
$ws.Range('E21').Value = '  +0.35%  '
$ws.Range('E22').Value = '  -2.30%  '
$ws.Range('D23').Value = "'79.87"
$ws.Range('E23').Value = '  -0.44%  '
$ws.Range('D24').Value = '3.832.44'
$ws.Range('E24').Value = '  +0.08%  '
$ws.Range('E25').Value = '  -0.09%  '
$ws.Range('E26').Value = '  -5.36%  '
$ws.Range('D27').Value = "'10.94"
$ws.Range('E27').Value = '  -4.36%  '
$ws.Range('D28').Value = "'9.12"
$ws.Range('E28').Value = '  -4.23%  '
$ws.Range('D29').Value = "'2.68"
$ws.Range('E29').Value = '  -1.41%  '
$ws.Range('E30').Value = '  -3.48%  '
$ws.Range('E31').Value = '  -3.22%  '
$ws.Range('E32').Value = '  +0.16%  '
$ws.Range('B33').Value = 'EthereumClassic'
$ws.Range('C33').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D33').Value = "'26.97"
$ws.Range('E33').Value = '  -0.80%  '
$ws.Range('B34').Value = 'ImmutableX'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D34').Value = "'1.99"
$ws.Range('E34').Value = '  -4.52%  '
$ws.Range('D35').Value = '3.674.66'
$ws.Range('E35').Value = '  +0.68%  '
$ws.Range('D36').Value = "'0.158"
$ws.Range('E36').Value = '  -5.14%  '
$ws.Range('D37').Value = "'8.25"
$ws.Range('E37').Value = '  -2.93%  '
$ws.Range('E38').Value = '  -1.86%  '
$ws.Range('D40').Value = "'2.24"
$ws.Range('E40').Value = '  -3.94%  '
$ws.Range('D42').Value = "'0.0905"
$ws.Range('E42').Value = '  -3.03%  '
$ws.Range('D43').Value = "'170.22"
$ws.Range('E43').Value = '  +4.63%  '
$ws.Range('D44').Value = "'0.942"
$ws.Range('E44').Value = '  -1.00%  '
$ws.Range('D45').Value = "'47.67"
$ws.Range('E45').Value = '  -0.65%  '
$ws.Range('D46').Value = "'0.000280"
$ws.Range('E46').Value = '  -2.47%  '
$ws.Range('B47').Value = 'dogwifhat'
$ws.Range('C47').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D47').Value = "'2.72"
$ws.Range('E47').Value = '  -4.40%  '
$ws.Range('B48').Value = 'InjectiveProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D48').Value = "'28.26"
$ws.Range('E48').Value = '  -5.28%  '
$ws.Range('E49').Value = '  -3.74%  '
$ws.Range('E50').Value = '  -5.46%  '
$ws.Range('D51').Value = "'7.79"
$ws.Range('E51').Value = '  -2.86%  '

# The quote-prefix entries above flip those cells to a "text" style;
# restore the default/Normal style so formatting matches the original file.
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D10').Style = 'Normal'
$ws.Range('D14').Style = 'Normal'
$ws.Range('D18').Style = 'Normal'
$ws.Range('D19').Style = 'Normal'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D23').Style = 'Normal'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D28').Style = 'Normal'
$ws.Range('D29').Style = 'Normal'
$ws.Range('D33').Style = 'Normal'
$ws.Range('D34').Style = 'Normal'
$ws.Range('D36').Style = 'Normal'
$ws.Range('D37').Style = 'Normal'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D42').Style = 'Normal'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D44').Style = 'Normal'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D46').Style = 'Normal'
$ws.Range('D47').Style = 'Normal'
$ws.Range('D48').Style = 'Normal'
$ws.Range('D51').Style = 'Normal'
